$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $newValue)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
}

Set-TextValue "D2" "29.482.06"
Set-TextValue "D3" "1.856.01"
Set-TextValue "E3" "  +1.20%  "
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "245.12"
Set-TextValue "E5" "  -0.42%  "
Set-TextValue "D6" "0.6939"
Set-TextValue "E6" "  +0.59%  "
Set-TextValue "D7" "1.0000"
Set-TextValue "E7" "  +0.01%  "
Set-TextValue "E8" "  +0.24%  "
Set-TextValue "D9" "0.07662"
Set-TextValue "E9" "  -0.48%  "
Set-TextValue "D10" "23.60"
Set-TextValue "E10" "  +0.13%  "
Set-TextValue "D11" "0.07771"
Set-TextValue "E11" "  -0.42%  "
Set-TextValue "D12" "5.137"
Set-TextValue "E12" "  +1.08%  "
Set-TextValue "D13" "1.844.79"
Set-TextValue "E13" "  +0.49%  "
Set-TextValue "B14" "Litecoin"
Set-TextValue "C14" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D14" "90.88"
Set-TextValue "E14" "  +0.43%  "
Set-TextValue "B15" "Polygon"
Set-TextValue "C15" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D15" "0.6907"
Set-TextValue "E15" "  +1.47%  "
Set-TextValue "D16" "6.351"
Set-TextValue "E16" "  -1.46%  "
Set-TextValue "D17" "29.475.78"
Set-TextValue "E17" "  +1.77%  "
Set-TextValue "D18" "0.000008273"
Set-TextValue "E18" "  -0.82%  "
Set-TextValue "D19" "2.101.47"
Set-TextValue "E19" "  +0.71%  "
Set-TextValue "D20" "237.64"
Set-TextValue "E20" "  -2.51%  "
Set-TextValue "D21" "12.71"
Set-TextValue "E21" "  -0.13%  "
Set-TextValue "D22" "0.9999"
Set-TextValue "E22" "  +0.04%  "
Set-TextValue "D23" "7.651"
Set-TextValue "E23" "  +2.34%  "
Set-TextValue "E24" "  +0.05%  "
Set-TextValue "D25" "0.1491"
Set-TextValue "E25" "  +1.35%  "
Set-TextValue "D26" "8.894"
Set-TextValue "E26" "  +1.00%  "
Set-TextValue "D27" "159.58"
Set-TextValue "E27" "  -1.95%  "
Set-TextValue "D28" "18.25"
Set-TextValue "E28" "  +0.21%  "
Set-TextValue "D29" "1.536"
Set-TextValue "E29" "  -1.26%  "
Set-TextValue "E30" "  +0.60%  "
Set-TextValue "D31" "4.155"
Set-TextValue "E31" "  -0.24%  "
Set-TextValue "D32" "1.199"
Set-TextValue "E32" "  +1.90%  "
Set-TextValue "E33" "  -0.42%  "
Set-TextValue "D34" "0.7720"
Set-TextValue "E34" "  +0.40%  "
Set-TextValue "E35" "  +2.34%  "
Set-TextValue "E36" "  +0.36%  "
Set-TextValue "E37" "  +0.17%  "
Set-TextValue "D38" "1.333.30"
Set-TextValue "E38" "  +7.64%  "
Set-TextValue "D39" "0.01865"
Set-TextValue "E39" "  +0.89%  "
Set-TextValue "D40" "2.715"
Set-TextValue "E40" "  +0.73%  "
Set-TextValue "D41" "0.9613"
Set-TextValue "E41" "  +3.48%  "
Set-TextValue "D42" "106.11"
Set-TextValue "E42" "  -2.35%  "
Set-TextValue "D43" "5.841"
Set-TextValue "E43" "  +0.44%  "
Set-TextValue "D44" "1.001"
Set-TextValue "E44" "  +0.14%  "
Set-TextValue "D45" "9.798"
Set-TextValue "E45" "  +2.28%  "
Set-TextValue "D46" "0.00000000124"
Set-TextValue "E46" "  +1.19%  "
Set-TextValue "D47" "1.999.95"
Set-TextValue "E47" "  +0.82%  "
Set-TextValue "D48" "0.5219"
Set-TextValue "E48" "  +0.94%  "
Set-TextValue "E49" "  +1.96%  "
Set-TextValue "E50" "  -2.04%  "
Set-TextValue "D51" "6.964"
Set-TextValue "E51" "  +0.43%  "

Write-Host "Cryptos list updated"
